$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "ECs", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 2, 1, 6.780879499999999, 13.561759, 0.09314755032665376, 0.07273600820493056, 2.491738145760666, 14.950428874564, 0.05190752325169092, 0.04053296113416156),
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "FAPs", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 3, 1, 24.18506433333333, 72.555193, 0.3322252662272683, 0.389136476570504, 8.887172726936443, 79.984554542428, 0.1851362775619506, 0.2168506915622518),
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "MuSCs", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 3, 1, 1.592065666666667, 4.776197, 0.02186987938794236, 0.02561625701948875, 0.5850289408902222, 5.265260468011999, 0.01218723701117514, 0.01427494821063397),
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "Neutrophils", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 2, 1, 25.158886, 50.31777200000001, 0.3456024546443489, 0.2698701456828592, 9.245018430285334, 55.470110581712, 0.1925908667204072, 0.1503881831872697),
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 3, 1, 14.00046133333333, 42.001384, 0.1923214645892228, 0.225266723235713, 5.144684190673778, 46.302157716064, 0.1071733057923238, 0.1255324228407979),
    @("Inflammatory-Mac", "Bmp10", "Bmpr1a", "Resolving-Mac", 1, 0.3333333333333333, 0.3674653333333333, 1.102396, 0.5572612813719676, 0.5572612813719677, 3, 1, 1.079828666666667, 3.239486, 0.01483338482456395, 0.01737438928650463, 0.3967996009395555, 3.571196408456, 0.008266071034420003, 0.009682074436852959),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "ECs", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 2, 1, 6.780879499999999, 13.561759, 0.09314755032665376, 0.07273600820493056, 1.979661947972833, 11.877971687837, 0.04124002707496283, 0.03220304707076901),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 3, 1, 24.18506433333333, 72.555193, 0.3322252662272683, 0.389136476570504, 7.060773100299889, 63.54695790269901, 0.1470889886653176, 0.1722857850082523),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "MuSCs", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 3, 1, 1.592065666666667, 4.776197, 0.02186987938794236, 0.02561625701948875, 0.4647998565634444, 4.183198709071, 0.009682642376767216, 0.01134130880885479),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "Neutrophils", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 2, 1, 25.158886, 50.31777200000001, 0.3456024546443489, 0.2698701456828592, 7.345078063632667, 44.070468381796, 0.1530115879239417, 0.1194819624955895),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "Inflammatory-Mac", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 3, 1, 14.00046133333333, 42.001384, 0.1923214645892228, 0.225266723235713, 4.087402018523556, 36.786618166712, 0.08514815879689899, 0.09973430039491518),
    @("Resolving-Mac", "Bmp10", "Bmpr1a", "Resolving-Mac", 2, 0.6666666666666666, 0.2919476666666667, 0.875843, 0.4427387186280323, 0.4427387186280324, 3, 1, 1.079828666666667, 3.239486, 0.01483338482456395, 0.01737438928650463, 0.3152534596331111, 2.837281136698, 0.006567313790143941, 0.007692314849651673)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowData = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

Write-Output "Updated $($data.Count) data rows with refreshed TPM values."
